# Insert a new data row at row 12 (pushing existing rows 12-51 down to 13-52)
# and populate it with the new observation described in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 12..51 down by one row, inheriting formatting from the row
# being split (matches Excel's default Insert behaviour).
$ws.Rows.Item(12).Insert()

# Fill in the new row 12 with the new record's values.
$ws.Range("A12").Value = 7
$ws.Range("B12").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C12").Value = "Ñuble"
$ws.Range("D12").Value = 44819
$ws.Range("E12").Value = 16
$ws.Range("F12").Value = 100112001
$ws.Range("G12").Value = "Berenjena"
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("I12").Value = "Segunda"
$ws.Range("J12").Value = 60
$ws.Range("K12").Value = 14000
$ws.Range("L12").Value = 14000
$ws.Range("M12").Value = 14000
$ws.Range("N12").Value = "`$/caja 90 unidades"
$ws.Range("O12").Value = "Región de Arica y Parinacota"
$ws.Range("P12").Value = 156
$ws.Range("Q12").Value = 90
$ws.Range("R12").Value = "Hortaliza"

# Keep the date formatting consistent with the rest of column D.
$ws.Range("D12").NumberFormat = "YYYY-MM-DD HH:MM:SS"
